$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to Text format so numeric-looking strings
# (e.g. "0.606", "5.10") are preserved exactly as text, matching the
# original inline-string cell contents instead of being auto-converted
# to floating point numbers by the smart-entry parser.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "37.017.94"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.994.73"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "243.97"
$ws.Range("E5").Value = "  -3.80%  "
$ws.Range("D6").Value = "0.606"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "54.69"
$ws.Range("E8").Value = "  -3.40%  "
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").Value = "57.25"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").Value = "0.0975"
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("D13").Value = "2.290.20"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "14.16"
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("D15").Value = "20.94"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "0.759"
$ws.Range("E16").Value = "  -6.22%  "
$ws.Range("D17").Value = "2.026.20"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "5.06"
$ws.Range("E18").Value = "  -4.45%  "
$ws.Range("D19").Value = "36.954.42"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "68.62"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "0.0₃0812"
$ws.Range("E21").Value = "  -3.47%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.10"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "228.85"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  -5.90%  "
$ws.Range("D26").Value = "2.35"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").Value = "162.58"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "8.70"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "19.25"
$ws.Range("E29").Value = "  -2.64%  "
$ws.Range("D30").Value = "0.127"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").Value = "1.30"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "4.44"
$ws.Range("E33").Value = "  -4.41%  "
$ws.Range("D34").Value = "0.0613"
$ws.Range("E34").Value = "  -6.49%  "
$ws.Range("D35").Value = "4.24"
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("E36").Value = "  -4.44%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "5.31"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("D42").Value = "1.434.37"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("D44").Value = "0.0205"
$ws.Range("E44").Value = "  -4.18%  "
$ws.Range("D45").Value = "0.0889"
$ws.Range("E45").Value = "  -7.33%  "
$ws.Range("D46").Value = "88.38"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "15.28"
$ws.Range("E47").Value = "  -4.15%  "
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "6.79"
$ws.Range("E50").Value = "  -6.78%  "
$ws.Range("D51").Value = "2.181.20"
$ws.Range("E51").Value = "  -1.46%  "
